$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3: corrected Hydrogen demand value for Iron & steel
$ws.Range("B3").Value = 1042154.906793731

# D3: Hydrogen / Non-metallic minerals value removed (now blank)
$ws.Range("D3").Value = ""

# C4: Methanol / Chemicals corrected to 0
$ws.Range("C4").Value = 0

# C5: Ammonia / Chemicals corrected value
$ws.Range("C5").Value = 1998.450691666861

# Row 7 label changed from "Other" to "Biogas", with corrected value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 3635.121303557358

# Insert a new row 8 ("Other"), copying the formatting from row 7
$ws.Range("A7:D7").Copy()
$ws.Range("A8:D8").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

$ws.Range("A8").Value = "Other"
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = 3938.178595147363
